$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps being stored as text (it holds values like
# "27.254.61" / "1.785.96" that are not valid numbers, and values like
# "1.002" that Excel would otherwise auto-convert to a number).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.254.61'
$ws.Range("E2").Value = '  -0.71%  '
$ws.Range("D3").Value = '1.785.96'
$ws.Range("E3").Value = '  -1.42%  '
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").Value = '335.42'
$ws.Range("E5").Value = '  -2.83%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.23%  '
$ws.Range("D7").Value = '0.3811'
$ws.Range("E7").Value = '  -0.72%  '
$ws.Range("D8").Value = '0.3415'
$ws.Range("E8").Value = '  -2.91%  '
$ws.Range("D9").Value = '48.35'
$ws.Range("E9").Value = '  -3.51%  '
$ws.Range("D10").Value = '1.194'
$ws.Range("E10").Value = '  -3.37%  '
$ws.Range("D11").Value = '0.07480'
$ws.Range("E11").Value = '  -3.71%  '
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("D13").Value = '21.86'
$ws.Range("E13").Value = '  -3.10%  '
$ws.Range("D14").Value = '6.453'
$ws.Range("E14").Value = '  -2.67%  '
$ws.Range("D15").Value = '1.787.92'
$ws.Range("E15").Value = '  -1.18%  '
$ws.Range("D16").Value = '7.077'
$ws.Range("E16").Value = '  -2.04%  '
$ws.Range("D17").Value = '0.00001095'
$ws.Range("E17").Value = '  -2.75%  '
$ws.Range("D18").Value = '0.06653'
$ws.Range("E18").Value = '  -1.60%  '
$ws.Range("D19").Value = '83.71'
$ws.Range("E19").Value = '  -3.75%  '
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("D21").Value = '6.623'
$ws.Range("E21").Value = '  +1.30%  '
$ws.Range("D22").Value = '17.31'
$ws.Range("E22").Value = '  -2.82%  '
$ws.Range("D23").Value = '27.245.08'
$ws.Range("E23").Value = '  -0.69%  '
$ws.Range("D24").Value = '12.33'
$ws.Range("E24").Value = '  -6.38%  '
$ws.Range("D25").Value = '2.411'
$ws.Range("E25").Value = '  -2.27%  '
$ws.Range("D26").Value = '1.488'
$ws.Range("E26").Value = '  -0.71%  '
$ws.Range("D27").Value = '2.537'
$ws.Range("E27").Value = '  -5.35%  '
$ws.Range("D28").Value = '21.26'
$ws.Range("E28").Value = '  -4.25%  '
$ws.Range("D29").Value = '153.48'
$ws.Range("E29").Value = '  -0.57%  '
$ws.Range("D30").Value = '1.989.03'
$ws.Range("E30").Value = '  -1.15%  '
$ws.Range("D31").Value = '134.03'
$ws.Range("E31").Value = '  -2.14%  '
$ws.Range("D32").Value = '4.015'
$ws.Range("E32").Value = '  -1.54%  '
$ws.Range("D33").Value = '6.065'
$ws.Range("E33").Value = '  -5.06%  '
$ws.Range("D34").Value = '0.08701'
$ws.Range("E34").Value = '  -1.32%  '
$ws.Range("D35").Value = '13.32'
$ws.Range("E35").Value = '  -4.24%  '
$ws.Range("D36").Value = '1.656'
$ws.Range("E36").Value = '  -3.83%  '
$ws.Range("D37").Value = '0.6934'
$ws.Range("E37").Value = '  -2.41%  '
$ws.Range("D38").Value = '5.438'
$ws.Range("E38").Value = '  -3.69%  '
$ws.Range("D39").Value = '0.2207'
$ws.Range("E39").Value = '  -2.75%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '0.06322'
$ws.Range("E40").Value = '  -3.21%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '8.787'
$ws.Range("E41").Value = '  -2.54%  '
$ws.Range("D42").Value = '0.02338'
$ws.Range("E42").Value = '  -3.54%  '
$ws.Range("E43").Value = '  -4.39%  '
$ws.Range("D44").Value = '14.43'
$ws.Range("E44").Value = '  -3.38%  '
$ws.Range("D45").Value = '0.6509'
$ws.Range("E45").Value = '  -1.57%  '
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("D47").Value = '3.848'
$ws.Range("E47").Value = '  -3.27%  '
$ws.Range("D48").Value = '2.144'
$ws.Range("E48").Value = '  -2.04%  '
$ws.Range("D49").Value = '128.87'
$ws.Range("E49").Value = '  -3.22%  '
$ws.Range("D50").Value = '0.07131'
$ws.Range("E50").Value = '  -3.20%  '
$ws.Range("D51").Value = '78.96'
$ws.Range("E51").Value = '  -2.18%  '
